$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Block 1 (rows 1-10): second timing series in columns F:G, mirroring
# the existing A:B series (n / time(ms)), using the same growth factor
# stored in $J$3.
# ---------------------------------------------------------------------

$ws.Range("F1").Value = "n"
$ws.Range("G1").Value = "time(ms)"

$ws.Range("F2").Value = 30
$ws.Range("G2").Value = 582

$ws.Range("F3").Formula = '=F2*$J$3'
$ws.Range("G3").Value = 1729

$ws.Range("F4:F9").Formula = '=F3*$J$3'
$ws.Range("G4").Value = 5305
$ws.Range("G5").Value = 16138

$ws.Range("B6").Value = "Oot"
$ws.Range("G6").Value = "Oot"

$ws.Range("G7:G9").ClearContents()

$ws.Range("A10").Value = "Substraction4"
$ws.Range("F10").Value = "Substraction5"

# Underline the one real data point of the second series (F2) to flag it.
$ws.Range("F2").Font.Underline = $true

# ---------------------------------------------------------------------
# Block 2 (rows 14-23): another pair of timing series (A:B and F:G)
# for a different experiment (division).
# ---------------------------------------------------------------------

$ws.Range("A14").Value = "n"
$ws.Range("B14").Value = "time(ms)"
$ws.Range("F14").Value = "n"
$ws.Range("G14").Value = "time(ms)"

$ws.Range("A15").Value = 1000
$ws.Range("B15").Value = 65
$ws.Range("F15").Value = 1000
$ws.Range("G15").Value = 35

$ws.Range("A16").Formula = '=A15*$J$3'
$ws.Range("B16").Value = 187
$ws.Range("F16").Formula = '=F15*$J$3'
$ws.Range("G16").Value = 134

$ws.Range("A17:A22").Formula = '=A16*$J$3'
$ws.Range("B17").Value = 739
$ws.Range("F17:F22").Formula = '=F16*$J$3'
$ws.Range("G17").Value = 529

$ws.Range("B18").Value = 2954
$ws.Range("G18").Value = 2994

$ws.Range("B19").Value = 11953
$ws.Range("G19").Value = 2096

$ws.Range("B20").Value = 47034
$ws.Range("G20").Value = 9108

$ws.Range("B21").Value = "Oot"
$ws.Range("G21").Value = 9108

$ws.Range("B22").Value = "Oot"
$ws.Range("G22").Value = "Oot"

$ws.Range("A23").Value = "division4"
$ws.Range("F23").Value = "Division5"
$ws.Range("B23").Value = "size = 5"
$ws.Range("G23").Value = "size = 1"

# Match the author's final selection.
$ws.Range("G22").Select() | Out-Null
